# Insert 9 new trading-day rows (2019-11-18 .. 2019-11-28) right before the
# existing 2019-11-29 row (currently row 346), shifting all subsequent rows
# down by 9 (dimension grows from A1:I418 to A1:I427).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows("346:354").Insert()

function Set-TextCell($addr, $text) {
    $rng = $ws.Range($addr)
    $rng.NumberFormat = "@"
    $rng.Value = $text
    $rng.Style = "Normal"
}

# --- row 346 : 2019-11-18 ---
$ws.Range("A346").Value = 1574035200
Set-TextCell "B346" "2019-11-18"
Set-TextCell "C346" "5286"
Set-TextCell "D346" "MI"
$ws.Range("E346").Value = 1.613
$ws.Range("F346").Value = 1.647
$ws.Range("G346").Value = 1.613
$ws.Range("H346").Value = 1.64
$ws.Range("I346").Value = 1601849

# --- row 347 : 2019-11-19 ---
$ws.Range("A347").Value = 1574121600
Set-TextCell "B347" "2019-11-19"
Set-TextCell "C347" "5286"
Set-TextCell "D347" "MI"
$ws.Range("E347").Value = 1.647
$ws.Range("F347").Value = 1.76
$ws.Range("G347").Value = 1.647
$ws.Range("H347").Value = 1.76
$ws.Range("I347").Value = 6734847

# --- row 348 : 2019-11-20 ---
$ws.Range("A348").Value = 1574208000
Set-TextCell "B348" "2019-11-20"
Set-TextCell "C348" "5286"
Set-TextCell "D348" "MI"
$ws.Range("E348").Value = 1.76
$ws.Range("F348").Value = 1.793
$ws.Range("G348").Value = 1.733
$ws.Range("H348").Value = 1.747
$ws.Range("I348").Value = 2616899

# --- row 349 : 2019-11-21 ---
$ws.Range("A349").Value = 1574294400
Set-TextCell "B349" "2019-11-21"
Set-TextCell "C349" "5286"
Set-TextCell "D349" "MI"
$ws.Range("E349").Value = 1.747
$ws.Range("F349").Value = 1.747
$ws.Range("G349").Value = 1.673
$ws.Range("H349").Value = 1.68
$ws.Range("I349").Value = 2618099

# --- row 350 : 2019-11-22 ---
$ws.Range("A350").Value = 1574380800
Set-TextCell "B350" "2019-11-22"
Set-TextCell "C350" "5286"
Set-TextCell "D350" "MI"
$ws.Range("E350").Value = 1.747
$ws.Range("F350").Value = 1.767
$ws.Range("G350").Value = 1.727
$ws.Range("H350").Value = 1.74
$ws.Range("I350").Value = 3625798

# --- row 351 : 2019-11-25 ---
$ws.Range("A351").Value = 1574640000
Set-TextCell "B351" "2019-11-25"
Set-TextCell "C351" "5286"
Set-TextCell "D351" "MI"
$ws.Range("E351").Value = 1.76
$ws.Range("F351").Value = 1.76
$ws.Range("G351").Value = 1.7
$ws.Range("H351").Value = 1.7
$ws.Range("I351").Value = 1341749

# --- row 352 : 2019-11-26 ---
$ws.Range("A352").Value = 1574726400
Set-TextCell "B352" "2019-11-26"
Set-TextCell "C352" "5286"
Set-TextCell "D352" "MI"
$ws.Range("E352").Value = 1.707
$ws.Range("F352").Value = 1.727
$ws.Range("G352").Value = 1.673
$ws.Range("H352").Value = 1.68
$ws.Range("I352").Value = 1255649

# --- row 353 : 2019-11-27 ---
$ws.Range("A353").Value = 1574812800
Set-TextCell "B353" "2019-11-27"
Set-TextCell "C353" "5286"
Set-TextCell "D353" "MI"
$ws.Range("E353").Value = 1.68
$ws.Range("F353").Value = 1.727
$ws.Range("G353").Value = 1.673
$ws.Range("H353").Value = 1.693
$ws.Range("I353").Value = 1583849

# --- row 354 : 2019-11-28 ---
$ws.Range("A354").Value = 1574899200
Set-TextCell "B354" "2019-11-28"
Set-TextCell "C354" "5286"
Set-TextCell "D354" "MI"
$ws.Range("E354").Value = 1.72
$ws.Range("F354").Value = 1.72
$ws.Range("G354").Value = 1.693
$ws.Range("H354").Value = 1.693
$ws.Range("I354").Value = 586650
